$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "88.395.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +8.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.324.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.53%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "650.82"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.394"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +33.98%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.322.28"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.588"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000282"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +8.60%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.10"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +9.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.941.67"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.45"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.159.61"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +8.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.327.30"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.59"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.15"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.65"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "454.87"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.46"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.39%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.53"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.72"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +12.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.508.11"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "78.26"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.81%  "

$ws.Range("B29").Value = "Cronos"
$ws.Range("C29").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.215"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +52.62%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000131"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.88%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.37"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "596.63"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.15%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.993"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.13"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.21"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +20.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.144"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.91%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.16"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.419"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.85"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.12"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.47"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.09%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "158.64"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "189.59"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.50"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.41"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.778"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.95%  "
